# Auto-generated COM-interop script to apply the W/X/Y/Z updates for rows 2-60,
# delete obsolete rows 61-68, and shrink the conditional-formatting / dimension range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update synthetics_from_this_seed (W), last_delta (X), last_neighbor_z (Y), timestamp (Z) ---

$ws.Cells.Item(2, 23).Value2 = 3
$ws.Cells.Item(2, 24).Value2 = 0.4538824667597043
$ws.Cells.Item(2, 25).Value2 = 30
$ws.Cells.Item(2, 26).Value2 = "2025-10-29T23:40:52.847678"

$ws.Cells.Item(3, 23).Value2 = 1
$ws.Cells.Item(3, 24).Value2 = 0.5616240759128834
$ws.Cells.Item(3, 25).Value2 = 123
$ws.Cells.Item(3, 26).Value2 = "2025-10-29T23:40:52.848854"

$ws.Cells.Item(4, 23).Value2 = 3
$ws.Cells.Item(4, 24).Value2 = 0.4329311706285884
$ws.Cells.Item(4, 25).Value2 = 24
$ws.Cells.Item(4, 26).Value2 = "2025-10-29T23:40:52.848854"

$ws.Cells.Item(5, 23).Value2 = 5
$ws.Cells.Item(5, 24).Value2 = 0.4153959819657586
$ws.Cells.Item(5, 25).Value2 = 6
$ws.Cells.Item(5, 26).Value2 = "2025-10-29T23:40:52.848854"

$ws.Cells.Item(6, 26).Value2 = "2025-10-29T23:40:52.848854"

$ws.Cells.Item(7, 26).Value2 = "2025-10-29T23:40:52.848854"

$ws.Cells.Item(8, 23).Value2 = 6
$ws.Cells.Item(8, 24).Value2 = 0.5079682182603347
$ws.Cells.Item(8, 25).Value2 = 11
$ws.Cells.Item(8, 26).Value2 = "2025-10-29T23:40:52.849855"

$ws.Cells.Item(9, 26).Value2 = "2025-10-29T23:40:52.849855"

$ws.Cells.Item(10, 26).Value2 = "2025-10-29T23:40:52.849855"

$ws.Cells.Item(11, 26).Value2 = "2025-10-29T23:40:52.849855"

$ws.Cells.Item(12, 23).Value2 = 4
$ws.Cells.Item(12, 24).Value2 = 0.5614880310328125
$ws.Cells.Item(12, 25).Value2 = 120
$ws.Cells.Item(12, 26).Value2 = "2025-10-29T23:40:52.849855"

$ws.Cells.Item(13, 26).Value2 = "2025-10-29T23:40:52.850852"

$ws.Cells.Item(14, 23).Value2 = 1
$ws.Cells.Item(14, 24).Value2 = 0.4062858371373469
$ws.Cells.Item(14, 25).Value2 = 61
$ws.Cells.Item(14, 26).Value2 = "2025-10-29T23:40:52.850852"

$ws.Cells.Item(15, 26).Value2 = "2025-10-29T23:40:52.850852"

$ws.Cells.Item(16, 23).Value2 = 4
$ws.Cells.Item(16, 24).Value2 = 0.405083825348819
$ws.Cells.Item(16, 25).Value2 = 21
$ws.Cells.Item(16, 26).Value2 = "2025-10-29T23:40:52.850852"

$ws.Cells.Item(17, 26).Value2 = "2025-10-29T23:40:52.850852"

$ws.Cells.Item(18, 26).Value2 = "2025-10-29T23:40:52.851852"

$ws.Cells.Item(19, 26).Value2 = "2025-10-29T23:40:52.851852"

$ws.Cells.Item(20, 26).Value2 = "2025-10-29T23:40:52.851852"

$ws.Cells.Item(21, 26).Value2 = "2025-10-29T23:40:52.851852"

$ws.Cells.Item(22, 26).Value2 = "2025-10-29T23:40:52.851852"

$ws.Cells.Item(23, 26).Value2 = "2025-10-29T23:40:52.852852"

$ws.Cells.Item(24, 26).Value2 = "2025-10-29T23:40:52.852852"

$ws.Cells.Item(25, 26).Value2 = "2025-10-29T23:40:52.852852"

$ws.Cells.Item(26, 23).Value2 = 3
$ws.Cells.Item(26, 24).Value2 = 0.4705137712668338
$ws.Cells.Item(26, 25).Value2 = 18
$ws.Cells.Item(26, 26).Value2 = "2025-10-29T23:40:52.852852"

$ws.Cells.Item(27, 23).Value2 = 3
$ws.Cells.Item(27, 24).Value2 = 0.5541934359909122
$ws.Cells.Item(27, 25).Value2 = 123
$ws.Cells.Item(27, 26).Value2 = "2025-10-29T23:40:52.852852"

$ws.Cells.Item(28, 23).Value2 = 2
$ws.Cells.Item(28, 24).Value2 = 0.4969659942717967
$ws.Cells.Item(28, 25).Value2 = 134
$ws.Cells.Item(28, 26).Value2 = "2025-10-29T23:40:52.853855"

$ws.Cells.Item(29, 23).Value2 = 4
$ws.Cells.Item(29, 24).Value2 = 0.4636006949943728
$ws.Cells.Item(29, 25).Value2 = 31
$ws.Cells.Item(29, 26).Value2 = "2025-10-29T23:40:52.853855"

$ws.Cells.Item(30, 23).Value2 = 9
$ws.Cells.Item(30, 24).Value2 = 0.4641560129943472
$ws.Cells.Item(30, 25).Value2 = 108
$ws.Cells.Item(30, 26).Value2 = "2025-10-29T23:40:52.881787"

$ws.Cells.Item(31, 23).Value2 = 7
$ws.Cells.Item(31, 24).Value2 = 0.4978905520555126
$ws.Cells.Item(31, 25).Value2 = 22
$ws.Cells.Item(31, 26).Value2 = "2025-10-29T23:40:52.882403"

$ws.Cells.Item(32, 23).Value2 = 9
$ws.Cells.Item(32, 24).Value2 = 0.5290345580818899
$ws.Cells.Item(32, 25).Value2 = 80
$ws.Cells.Item(32, 26).Value2 = "2025-10-29T23:40:52.882403"

$ws.Cells.Item(33, 23).Value2 = 7
$ws.Cells.Item(33, 24).Value2 = 0.5381875476204931
$ws.Cells.Item(33, 25).Value2 = 6
$ws.Cells.Item(33, 26).Value2 = "2025-10-29T23:40:52.882935"

$ws.Cells.Item(34, 23).Value2 = 5
$ws.Cells.Item(34, 24).Value2 = 0.5966846281789686
$ws.Cells.Item(34, 25).Value2 = 91
$ws.Cells.Item(34, 26).Value2 = "2025-10-29T23:40:52.882935"

$ws.Cells.Item(35, 26).Value2 = "2025-10-29T23:40:52.882935"

$ws.Cells.Item(36, 23).Value2 = 9
$ws.Cells.Item(36, 24).Value2 = 0.4081550283109528
$ws.Cells.Item(36, 25).Value2 = 135
$ws.Cells.Item(36, 26).Value2 = "2025-10-29T23:40:52.882935"

$ws.Cells.Item(37, 26).Value2 = "2025-10-29T23:40:52.882935"

$ws.Cells.Item(38, 23).Value2 = 13
$ws.Cells.Item(38, 24).Value2 = 0.4739308912122809
$ws.Cells.Item(38, 25).Value2 = 11
$ws.Cells.Item(38, 26).Value2 = "2025-10-29T23:40:52.883932"

$ws.Cells.Item(39, 26).Value2 = "2025-10-29T23:40:52.883932"

$ws.Cells.Item(40, 26).Value2 = "2025-10-29T23:40:52.883932"

$ws.Cells.Item(41, 23).Value2 = 7
$ws.Cells.Item(41, 24).Value2 = 0.4727259204758588
$ws.Cells.Item(41, 25).Value2 = 118
$ws.Cells.Item(41, 26).Value2 = "2025-10-29T23:40:52.883932"

$ws.Cells.Item(42, 26).Value2 = "2025-10-29T23:40:52.883932"

$ws.Cells.Item(43, 26).Value2 = "2025-10-29T23:40:52.927950"

$ws.Cells.Item(44, 23).Value2 = 4
$ws.Cells.Item(44, 24).Value2 = 0.4822074026636463
$ws.Cells.Item(44, 25).Value2 = 16
$ws.Cells.Item(44, 26).Value2 = "2025-10-29T23:40:52.927950"

$ws.Cells.Item(45, 23).Value2 = 4
$ws.Cells.Item(45, 24).Value2 = 0.42961738599068
$ws.Cells.Item(45, 25).Value2 = 116
$ws.Cells.Item(45, 26).Value2 = "2025-10-29T23:40:52.928950"

$ws.Cells.Item(46, 23).Value2 = 5
$ws.Cells.Item(46, 24).Value2 = 0.5721461166512687
$ws.Cells.Item(46, 25).Value2 = 16
$ws.Cells.Item(46, 26).Value2 = "2025-10-29T23:40:52.928950"

$ws.Cells.Item(47, 26).Value2 = "2025-10-29T23:40:52.928950"

$ws.Cells.Item(48, 26).Value2 = "2025-10-29T23:40:52.928950"

$ws.Cells.Item(49, 23).Value2 = 9
$ws.Cells.Item(49, 24).Value2 = 0.4690142496053366
$ws.Cells.Item(49, 25).Value2 = 115
$ws.Cells.Item(49, 26).Value2 = "2025-10-29T23:40:52.930138"

$ws.Cells.Item(50, 23).Value2 = 4
$ws.Cells.Item(50, 24).Value2 = 0.453356202855057
$ws.Cells.Item(50, 25).Value2 = 52
$ws.Cells.Item(50, 26).Value2 = "2025-10-29T23:40:52.930138"

$ws.Cells.Item(51, 23).Value2 = 6
$ws.Cells.Item(51, 24).Value2 = 0.4444215620941461
$ws.Cells.Item(51, 25).Value2 = 17
$ws.Cells.Item(51, 26).Value2 = "2025-10-29T23:40:52.930138"

$ws.Cells.Item(52, 23).Value2 = 5
$ws.Cells.Item(52, 24).Value2 = 0.5792182599846987
$ws.Cells.Item(52, 25).Value2 = 44
$ws.Cells.Item(52, 26).Value2 = "2025-10-29T23:40:52.930946"

$ws.Cells.Item(53, 26).Value2 = "2025-10-29T23:40:52.930946"

$ws.Cells.Item(54, 26).Value2 = "2025-10-29T23:40:52.930946"

$ws.Cells.Item(55, 26).Value2 = "2025-10-29T23:40:52.930946"

$ws.Cells.Item(56, 26).Value2 = "2025-10-29T23:40:52.930946"

$ws.Cells.Item(57, 23).Value2 = 12
$ws.Cells.Item(57, 24).Value2 = 0.5392608545679577
$ws.Cells.Item(57, 25).Value2 = 118
$ws.Cells.Item(57, 26).Value2 = "2025-10-29T23:40:52.931947"

$ws.Cells.Item(58, 26).Value2 = "2025-10-29T23:40:52.931947"

$ws.Cells.Item(59, 23).Value2 = 9
$ws.Cells.Item(59, 24).Value2 = 0.450783082786869
$ws.Cells.Item(59, 25).Value2 = 144
$ws.Cells.Item(59, 26).Value2 = "2025-10-29T23:40:52.931947"

$ws.Cells.Item(60, 23).Value2 = 6
$ws.Cells.Item(60, 24).Value2 = 0.4646405864041511
$ws.Cells.Item(60, 25).Value2 = 22
$ws.Cells.Item(60, 26).Value2 = "2025-10-29T23:40:52.931947"

# --- Remove the now-obsolete sample rows 61-68 ---
$ws.Rows("61:68").Delete()

# --- Shrink the conditional formatting range to match the new data extent ---
$fcs = $ws.Range("A2:Z68").FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:Z60"))
